$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Updating Item State:" heading -> append a new bold run:
#       " - Not sure on this yet"
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Updating Item State:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Collapse(0)   # wdCollapseEnd - put insertion point right after "State:"
    $rng.InsertAfter(" - Not sure on this yet")
    $rng.Font.Bold = $true
    $rng.Font.Size = 12
}

# ---------------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark from the end of the last
#    "Some text about what the thing does and the diagram below" paragraph
#    to wrap around the final diagram picture/paragraph instead.
#    Re-adding a bookmark under the same name relocates it (Word keeps
#    bookmark names unique), so this both removes the old one and creates
#    the new one in a single call.
# ---------------------------------------------------------------------------
$lastParaIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastParaIndex)
$d.Bookmarks.Add("_GoBack", $lastPara.Range)
